# Corrects quantities (Qty, col F) and their dependent Value figures (col G) for a
# batch of line items across the stock report, fixes the swapped SIG-3W / SIG-3w
# 'Lilliput LED Torch & Table Lamp' rows (11-12, all of cols B-G), and updates every
# Sub Total / Grand Total row (col B) downstream of those line items so the report
# stays internally consistent.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- ALPHA TRADING AND PROMOTIONS --
$ws.Range("B11").Value = 47438
$ws.Range("C11").Value = "SIG-3w Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D11").Value = 401.81
$ws.Range("E11").Value = 480.05
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 803.62
$ws.Range("B12").Value = 59408
$ws.Range("C12").Value = "SIG-3W Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D12").Value = 388.17
$ws.Range("E12").Value = 463.78
$ws.Range("F12").Value = 27
$ws.Range("G12").Value = 10480.59
# -- BHAWAR SALES CORPORATION --
$ws.Range("F60").Value = 20
$ws.Range("G60").Value = 2508
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("F78").Value = 10
$ws.Range("G78").Value = 940.9
# Row 83: Sub Total for BHAWAR SALES CORPORATION
$ws.Range("B83").Value = 104011.02
# -- Cholayil Pvt Ltd --
$ws.Range("F100").Value = 46
$ws.Range("G100").Value = 2276.08
$ws.Range("F102").Value = 100
$ws.Range("G102").Value = 4948
$ws.Range("F105").Value = 19
$ws.Range("G105").Value = 804.08
# Row 112: Sub Total for Cholayil Pvt Ltd
$ws.Range("B112").Value = 60656.05
# -- DABUR INDIA LIMITED --
$ws.Range("F146").Value = 7
$ws.Range("G146").Value = 357.14
# Row 160: Sub Total for DABUR INDIA LIMITED
$ws.Range("B160").Value = 13709.37
# -- EMAMI LTD --
$ws.Range("F173").Value = 55
$ws.Range("G173").Value = 4455
$ws.Range("F174").Value = 18
$ws.Range("G174").Value = 1845.9
# Row 175: Sub Total for EMAMI LTD
$ws.Range("B175").Value = 22685.68
# -- Glaxosmithkline Asia Private Limited --
$ws.Range("F189").Value = 121
$ws.Range("G189").Value = 8678.12
# Row 194: Sub Total for Glaxosmithkline Asia Private Limited
$ws.Range("B194").Value = 13334.57
# -- GODREJ CONSUMER PRODUCTS LIMITED --
$ws.Range("F196").Value = 12
$ws.Range("G196").Value = 1403.04
$ws.Range("F197").Value = 852
$ws.Range("G197").Value = 15762
# Row 204: Sub Total for GODREJ CONSUMER PRODUCTS LIMITED
$ws.Range("B204").Value = 21802.49
# -- HIMALAYA WELLNESS COMPANY --
$ws.Range("F236").Value = 28
$ws.Range("G236").Value = 897.96
$ws.Range("F247").Value = 129
$ws.Range("G247").Value = 5966.25
$ws.Range("F253").Value = 31
$ws.Range("G253").Value = 2542.93
$ws.Range("F256").Value = 18
$ws.Range("G256").Value = 910.08
$ws.Range("F257").Value = 269
$ws.Range("G257").Value = 8957.7
$ws.Range("F264").Value = 20
$ws.Range("G264").Value = 641.4
$ws.Range("F265").Value = 71
$ws.Range("G265").Value = 2276.97
$ws.Range("F268").Value = 91
$ws.Range("G268").Value = 2446.99
# Row 276: Sub Total for HIMALAYA WELLNESS COMPANY
$ws.Range("B276").Value = 251410.72
# -- HINDUSTAN UNILIVER LTD --
$ws.Range("F278").Value = 36
$ws.Range("G278").Value = 275.04
$ws.Range("F279").Value = 111
$ws.Range("G279").Value = 19209.66
$ws.Range("F280").Value = 3
$ws.Range("G280").Value = 917.52
$ws.Range("F283").Value = 50
$ws.Range("G283").Value = 4147
$ws.Range("F286").Value = 276
$ws.Range("G286").Value = 47690.04
$ws.Range("F287").Value = 22
$ws.Range("G287").Value = 3801.38
$ws.Range("F290").Value = 99
$ws.Range("G290").Value = 8133.84
$ws.Range("F293").Value = 32
$ws.Range("G293").Value = 5148.8
$ws.Range("F295").Value = 336
$ws.Range("G295").Value = 24766.56
$ws.Range("F296").Value = 58
$ws.Range("G296").Value = 9202.86
$ws.Range("F305").Value = 174
$ws.Range("G305").Value = 24975.96
$ws.Range("F306").Value = 86
$ws.Range("G306").Value = 9704.24
$ws.Range("F316").Value = 104
$ws.Range("G316").Value = 1521.52
$ws.Range("F326").Value = 251
$ws.Range("G326").Value = 27564.82
$ws.Range("F334").Value = 198
$ws.Range("G334").Value = 25118.28
$ws.Range("F335").Value = 96
$ws.Range("G335").Value = 9427.2
$ws.Range("F336").Value = 113
$ws.Range("G336").Value = 6500.89
$ws.Range("F337").Value = 89
$ws.Range("G337").Value = 7629.08
$ws.Range("F338").Value = 56
$ws.Range("G338").Value = 3861.76
$ws.Range("F349").Value = 88
$ws.Range("G349").Value = 19063.44
$ws.Range("F354").Value = 27
$ws.Range("G354").Value = 2805.57
$ws.Range("F356").Value = 77
$ws.Range("G356").Value = 4553.01
$ws.Range("F359").Value = 391
$ws.Range("G359").Value = 8226.64
$ws.Range("F360").Value = 398
$ws.Range("G360").Value = 33073.8
$ws.Range("F362").Value = 844
$ws.Range("G362").Value = 144602.52
$ws.Range("F363").Value = 392
$ws.Range("G363").Value = 59258.64
$ws.Range("F364").Value = 48
$ws.Range("G364").Value = 20087.52
$ws.Range("F365").Value = 2
$ws.Range("G365").Value = 159.82
$ws.Range("F366").Value = 128
$ws.Range("G366").Value = 20520.96
$ws.Range("F375").Value = 525
$ws.Range("G375").Value = 21619.5
$ws.Range("F377").Value = 571
$ws.Range("G377").Value = 22697.25
$ws.Range("F378").Value = 15
$ws.Range("G378").Value = 4305.15
$ws.Range("F379").Value = 205
$ws.Range("G379").Value = 29450.3
# Row 380: Sub Total for HINDUSTAN UNILIVER LTD
$ws.Range("B380").Value = 1105237.66
# -- HINDUSTAN UNILIVER LTD(GSK) --
$ws.Range("F382").Value = 175
$ws.Range("G382").Value = 32126.5
$ws.Range("F383").Value = 165
$ws.Range("G383").Value = 30290.7
$ws.Range("F389").Value = 41
$ws.Range("G389").Value = 8540.3
$ws.Range("F394").Value = 26
$ws.Range("G394").Value = 4219.54
# Row 398: Sub Total for HINDUSTAN UNILIVER LTD(GSK)
$ws.Range("B398").Value = 150538.18
# -- JNTL Consumer Health (India) Private Limited --
$ws.Range("F419").Value = 72
$ws.Range("G419").Value = 5499.36
# Row 422: Sub Total for JNTL Consumer Health (India) Private Limited
$ws.Range("B422").Value = 45098.2
# -- LIFE STYLE FOODS PVT LTD --
$ws.Range("F478").Value = 60
$ws.Range("G478").Value = 6172.2
# Row 491: Sub Total for LIFE STYLE FOODS PVT LTD
$ws.Range("B491").Value = 59547.31
# -- NETWAY HOME PRODUCTS INDIA PVT LTD --
$ws.Range("F546").Value = 28
$ws.Range("G546").Value = 616.56
# Row 548: Sub Total for NETWAY HOME PRODUCTS INDIA PVT LTD
$ws.Range("B548").Value = 3136.13
# -- RECKITT BENCKISER INDIA PVT LTD --
$ws.Range("F588").Value = 8
$ws.Range("G588").Value = 572.8
# Row 601: Sub Total for RECKITT BENCKISER INDIA PVT LTD
$ws.Range("B601").Value = 27774
# -- SOUTHERN HEALTH FOODS PVT LTD --
$ws.Range("F663").Value = 91
$ws.Range("G663").Value = 6511.96
$ws.Range("F665").Value = 173
$ws.Range("G665").Value = 15872.75
$ws.Range("F666").Value = 171
$ws.Range("G666").Value = 14620.5
$ws.Range("F670").Value = 31
$ws.Range("G670").Value = 4429.59
# Row 673: Sub Total for SOUTHERN HEALTH FOODS PVT LTD
$ws.Range("B673").Value = 86957.79
# -- TATA CONSUMER PRODUCT LIMITED --
$ws.Range("F709").Value = 459
$ws.Range("G709").Value = 37436.04
$ws.Range("F712").Value = 224
$ws.Range("G712").Value = 34612.48
$ws.Range("F713").Value = 261
$ws.Range("G713").Value = 21287.16
$ws.Range("F714").Value = 507
$ws.Range("G714").Value = 67481.7
$ws.Range("F718").Value = 256
$ws.Range("G718").Value = 5560.32
$ws.Range("F729").Value = 825
$ws.Range("G729").Value = 118800
$ws.Range("F731").Value = 649
$ws.Range("G731").Value = 78340.79
# Row 733: Sub Total for TATA CONSUMER PRODUCT LIMITED
$ws.Range("B733").Value = 646061.5
# -- VVD AND SONS PRIVATE LIMITED --
$ws.Range("F774").Value = 7
$ws.Range("G774").Value = 662.2
$ws.Range("F776").Value = 125
$ws.Range("G776").Value = 3261.25
$ws.Range("F778").Value = 2683
$ws.Range("G778").Value = 437624.13
# Row 783: Sub Total for VVD AND SONS PRIVATE LIMITED
$ws.Range("B783").Value = 451685.68
# Row 795: Grand Sub Total (sum of every company Sub Total)
$ws.Range("B795").Value = 5851461.17
# Row 796: Grand Total
$ws.Range("B796").Value = 5851461.17
